# material expense recording completed. disabling in progress.
# adding rounding to quantity and amount (#.########)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Notes \ Field" explanations in column B ---------------------
# New note explaining why the unauthorized-salesclerk (user18) test cases fail.
$wrongValues = "It fails due to wrong values since user 18 is a salesclerk. Hence, user 18 can not retrieve  providers, materials, so can't populate the corresponding combos."
$negWrongUnit = "It fails due to negative values and wrong unit."
$wrongUnit = "It fails due to wrong unit."

# Rows that previously said "It fails due to not enough rights." and now
# explain the real (wrong-values) failure reason; B15/B41/B45/B47 additionally
# get their font re-applied (same Arial 10, just "touched" through a style).
$ws.Range("B15").Value = $wrongValues
$ws.Range("B15").Style = "Normal"

$ws.Range("B20").Value = $wrongValues
$ws.Range("B27").Value = $wrongValues
$ws.Range("B31").Value = $negWrongUnit
$ws.Range("B33").Value = $wrongUnit
$ws.Range("B34").Value = $wrongValues
$ws.Range("B37").Value = $wrongUnit

$ws.Range("B41").Value = $wrongValues
$ws.Range("B41").Style = "Normal"

$ws.Range("B42").Value = $wrongUnit
$ws.Range("B43").Value = $wrongUnit

$ws.Range("B45").Value = $wrongValues
$ws.Range("B45").Style = "Normal"

$ws.Range("B46").Value = $wrongUnit

$ws.Range("B47").Value = $wrongValues
$ws.Range("B47").Style = "Normal"

# --- New test case row 49: checking recorded expenses --------------------
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "check recorded expenses"

# --- Leave the view focused on the newly added row ------------------------
$ws.Range("A50").Select()
